# Add a new "Лист8" worksheet (Test Case for getGrayCode) by cloning the
# existing last sheet ("Лист7") and then updating its text, sizes and a
# couple of border tweaks so it matches the new Test Case layout.

$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item(7)

# --- Create the new sheet right after Лист7 ------------------------------
$ws7.Copy([System.Reflection.Missing]::Value, $ws7)
$ws8 = $wb.Worksheets.Item(8)
$ws8.Name = "Лист8"

# --- Column widths on the new sheet ---------------------------------------
$ws8.Columns.Item(2).ColumnWidth = 18.5
$ws8.Columns.Item(3).ColumnWidth = 17.333333333333336
$ws8.Columns.Item(4).ColumnWidth = 16.833333333333336
$ws8.Columns.Item(5).ColumnWidth = 11.166666666666668
$ws8.Columns.Item(6).ColumnWidth = 10.666666666666668

# --- Row heights that changed on the new sheet ----------------------------
$ws8.Rows.Item(3).RowHeight = 23.25
$ws8.Rows.Item(4).RowHeight = 21.75
$ws8.Rows.Item(12).RowHeight = 30.75
$ws8.Rows.Item(15).RowHeight = 75

# --- Drop the old "Completed with unit testing" note row & the stray J
#     column helper cells that the new Test Case sheet doesn't carry -------
$ws8.Range("B5:C5").Clear()
$ws8.Range("J9:J10").Clear()

# --- Make row 3's merged value cells vertically centred like row 4 --------
$ws8.Range("C3:F3").VerticalAlignment = -4108

# --- Remove the inner vertical divider between columns C and D on the
#     header row so the whole row uses a uniform border -------------------
$ws8.Range("C9").Borders.Item(10).LineStyle = 1
$ws8.Range("C9").Borders.Item(10).Weight = -4138
$ws8.Range("C9").Borders.Item(10).Color = 0
$ws8.Range("D9").Borders.Item(7).LineStyle = -4142

# --- New Test Case text for getGrayCode (keep this insertion order so new
#     shared-string ids come out in the same sequence as the source edit) --
$ws8.Range("C3").Value2 = "getGrayCode"
$ws8.Range("A10").Value2 = "TC8"
$ws8.Range("B10").Value2 = "This test case checks the functionality of the function getGrayCode"
$ws8.Range("C10").Value2 = "Check if the function takes the input data"
$ws8.Range("C13").Value2 = "Check if the function successfully converts the input into uniqe gray code"
$ws8.Range("D10").Value2 = "The function takes the given input, and converts it into uniqe grey code"

# --- Selections: Лист7 keeps D10:D15 selected, Лист8 becomes the active
#     tab with I16 selected -------------------------------------------------
$ws7.Activate()
$ws7.Range("D10:D15").Select()

$ws8.Activate()
$ws8.Range("I16").Select()
